$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin ranking values (price/volume refresh + row shift for new GateToken entry,
# plus a BOLO/CoinbaseStockToken swap) captured by the Jan 29 2023 GitHub Actions data pull.
$updates = @(
    @{ Cell = 'D2'; Value = '308.35' }
    @{ Cell = 'E2'; Value = '-0.44%' }
    @{ Cell = 'D3'; Value = '39.43' }
    @{ Cell = 'E3'; Value = '0.79%' }
    @{ Cell = 'D4'; Value = '5.142' }
    @{ Cell = 'E4'; Value = '0.79%' }
    @{ Cell = 'D5'; Value = '0.08121' }
    @{ Cell = 'E5'; Value = '-0.31%' }
    @{ Cell = 'D6'; Value = '1.943' }
    @{ Cell = 'E6'; Value = '-2.18%' }
    @{ Cell = 'B7'; Value = 'GateToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D7'; Value = '4.240' }
    @{ Cell = 'E7'; Value = '0.92%' }
    @{ Cell = 'B8'; Value = 'KuCoinToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' }
    @{ Cell = 'D8'; Value = '8.139' }
    @{ Cell = 'E8'; Value = '2.66%' }
    @{ Cell = 'B9'; Value = 'MXToken' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D9'; Value = '0.9270' }
    @{ Cell = 'E9'; Value = '-0.57%' }
    @{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D10'; Value = '0.1414' }
    @{ Cell = 'E10'; Value = '-0.43%' }
    @{ Cell = 'B11'; Value = 'WazirX' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D11'; Value = '0.1929' }
    @{ Cell = 'E11'; Value = '-1.05%' }
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D12'; Value = '0.09028' }
    @{ Cell = 'E12'; Value = '-1.15%' }
    @{ Cell = 'B13'; Value = 'BitrueCoin' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D13'; Value = '0.03516' }
    @{ Cell = 'E13'; Value = '0.53%' }
    @{ Cell = 'B14'; Value = 'BitMartToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D14'; Value = '0.09813' }
    @{ Cell = 'E14'; Value = '-0.21%' }
    @{ Cell = 'B15'; Value = 'BitForexToken' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D15'; Value = '0.001403' }
    @{ Cell = 'E15'; Value = '-0.81%' }
    @{ Cell = 'B16'; Value = 'TigerCash' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D16'; Value = '0.005945' }
    @{ Cell = 'E16'; Value = '-1.19%' }
    @{ Cell = 'B17'; Value = 'LEO' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D17'; Value = '3.909' }
    @{ Cell = 'E17'; Value = '9.34%' }
    @{ Cell = 'D18'; Value = '3.375' }
    @{ Cell = 'E18'; Value = '-3.26%' }
    @{ Cell = 'E19'; Value = '0.14%' }
    @{ Cell = 'D20'; Value = '0.1313' }
    @{ Cell = 'E20'; Value = '-0.81%' }
    @{ Cell = 'D21'; Value = '4.758' }
    @{ Cell = 'E21'; Value = '-1.43%' }
    @{ Cell = 'D22'; Value = '0.2428' }
    @{ Cell = 'E22'; Value = '-1.61%' }
    @{ Cell = 'D23'; Value = '0.04384' }
    @{ Cell = 'E23'; Value = '-1.74%' }
    @{ Cell = 'D24'; Value = '0.001231' }
    @{ Cell = 'E24'; Value = '-0.84%' }
    @{ Cell = 'D25'; Value = '0.004843' }
    @{ Cell = 'E25'; Value = '-0.54%' }
    @{ Cell = 'E26'; Value = '-0.01%' }
    @{ Cell = 'D27'; Value = '0.0004006' }
    @{ Cell = 'E27'; Value = '-9.92%' }
    @{ Cell = 'D39'; Value = '0.02065' }
    @{ Cell = 'E39'; Value = '-2.86%' }
    @{ Cell = 'D40'; Value = '0.05094' }
    @{ Cell = 'E40'; Value = '-1.00%' }
    @{ Cell = 'D41'; Value = '0.007424' }
    @{ Cell = 'E41'; Value = '-0.73%' }
    @{ Cell = 'D42'; Value = '0.009838' }
    @{ Cell = 'E42'; Value = '-1.75%' }
    @{ Cell = 'D43'; Value = '0.1362' }
    @{ Cell = 'E43'; Value = '-0.21%' }
    @{ Cell = 'E44'; Value = '-0.01%' }
    @{ Cell = 'D45'; Value = '0.008494' }
    @{ Cell = 'E45'; Value = '-12.87%' }
    @{ Cell = 'D46'; Value = '0.00006409' }
    @{ Cell = 'E46'; Value = '3.22%' }
    @{ Cell = 'E47'; Value = '-0.02%' }
    @{ Cell = 'B48'; Value = 'CoinbaseStockToken' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin' }
    @{ Cell = 'D48'; Value = '0.001301' }
    @{ Cell = 'E48'; Value = '-18.86%' }
    @{ Cell = 'B49'; Value = 'BOLO' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo' }
    @{ Cell = 'D49'; Value = '0.002573' }
    @{ Cell = 'E49'; Value = '-100.00%' }
    @{ Cell = 'D50'; Value = '0.00002103' }
    @{ Cell = 'E50'; Value = '-0.02%' }
    @{ Cell = 'D51'; Value = '0.0002003' }
    @{ Cell = 'E51'; Value = '-0.02%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text interpretation so numeric-looking strings (prices, percentages)
    # keep their exact original formatting (leading/trailing zeros, percent sign, etc.)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    # Revert to the default (unstyled) cell style so we don't leave a stray number format behind
    $rng.Style = "Normal"
}
